# Apply the journal/planning updates recorded in the commit.
$wb = $excel.ActiveWorkbook

$wsPlanning = $wb.Worksheets.Item("Planning")
$wsJournal  = $wb.Worksheets.Item("Journal de travail")

# --- "Planning" sheet -------------------------------------------------
# Mark the "API Gateway" functioning-test column as done.
$wsPlanning.Range("G22").Value = "X"

# --- "Journal de travail" sheet ---------------------------------------
# Correct a previously logged entry's duration.
$wsJournal.Range("C13").Value = 1.5

# Log new work entries.
$wsJournal.Range("A14").Value = 45016
$wsJournal.Range("B14").Value = "mise en place readME sur gitHub"
$wsJournal.Range("C14").Value = 0.25

$wsJournal.Range("A15").Value = 45016
$wsJournal.Range("B15").Value = "commencement de l'inplémentation sur les REST"
$wsJournal.Range("C15").Value = 1.75

$wsJournal.Range("A16").Value = 45016
$wsJournal.Range("A17").Value = 45016

# --- Selection / active sheet ------------------------------------------
$wsJournal.Range("B16").Select()
$wsPlanning.Range("O18").Select()
$wsPlanning.Activate()
